$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$newDate = (Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0 -Millisecond 0).AddDays(45183)

foreach ($row in 2..7) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
